$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates ---
# B1: "dbExcel" -> "StatQuery"
$ws.Range("B1").Value = "StatQuery"

# B2: old Neo4jData filename -> new StatQuery Cypher text, with wrap-text style matching A2
$statQuery = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Endometrioid endometrial adenocarcinoma'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"
$ws.Range("B2").Value = $statQuery

# Give B2 the same wrap-text style already used by A2
$ws.Range("B2").WrapText = $true

# --- Column widths: column B widens to match column A (75.81640625 chars) ---
# Column A already has the correct width; leave it untouched so it keeps its
# exact stored value. Column B only needs to grow to match it.
$ws.Range("B1").EntireColumn.ColumnWidth = 75

# --- Selection moves to A4 ---
$ws.Range("A4").Select() | Out-Null

# --- Window view size/position ---
$win = $wb.Windows.Item(1)
$win.Left = -110
$win.Top = -110
$win.Width = 19420
$win.Height = 10420
